$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update Naufal's task status (row 2, column F) from WAITING to DONE
$ws.Range("F2").Value = "DONE"

# Move the active selection to F3, matching the post-edit cursor position
$ws.Range("F3").Select()
